# Update the "FilesTab" Neo4j query in cell B4: remove the `File Type` and
# `Breed` columns from the RETURN clause (ICDC Breed script correction).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newFileQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Flat-Coated Retriever'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Cells.Item(4, 2).Value2 = $newFileQuery

# The shorter query text now needs two fewer wrapped lines of row height.
$ws.Rows.Item(4).RowHeight = 217.5

# Reflect that B4 was the last-edited / selected cell, matching the author's
# workbook view state after making the change.
$ws.Cells.Item(4, 2).Select() | Out-Null
